$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71-169 down to 72-170.
$ws.Range("A71").EntireRow.Insert()

# Populate the newly inserted row 71 with the new record's data.
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 44915
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112043
$ws.Range("G71").Value = "Pepino ensalada"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 100
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 16000
$ws.Range("M71").Value = 15500
$ws.Range("N71").Value = "$/caja 50 unidades"
$ws.Range("O71").Value = "Región de Arica y Parinacota"
$ws.Range("P71").Value = 310
$ws.Range("Q71").Value = 50
$ws.Range("R71").Value = "Hortaliza"
